# Update LR-pair (Mfge8-Itgb3) NATMI output values following Dr Hou's advice:
# ligand/receptor-expressing cell counts (E,K) change from 1 to 3 for every data row,
# which changes the corresponding average/total expression values and specificity scores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.21340333333333
$ws.Range("H2").Value = 42.64021
$ws.Range("I2").Value = 0.07497543485230342
$ws.Range("J2").Value = 0.07497543485230343
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 102.5370645966355
$ws.Range("R2").Value = 922.8335813697199
$ws.Range("S2").Value = 0.03515466030495176
$ws.Range("T2").Value = 0.03515466030495177

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.21340333333333
$ws.Range("H3").Value = 42.64021
$ws.Range("I3").Value = 0.07497543485230342
$ws.Range("J3").Value = 0.07497543485230343
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 101.0622818667689
$ws.Range("R3").Value = 909.5605368009199
$ws.Range("S3").Value = 0.03464903352407965
$ws.Range("T3").Value = 0.03464903352407966

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.21340333333333
$ws.Range("H4").Value = 42.64021
$ws.Range("I4").Value = 0.07497543485230342
$ws.Range("J4").Value = 0.07497543485230343
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 15.08463284185444
$ws.Range("R4").Value = 135.76169557669
$ws.Range("S4").Value = 0.005171741023272009
$ws.Range("T4").Value = 0.00517174102327201

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.428335
$ws.Range("H5").Value = 88.285005
$ws.Range("I5").Value = 0.1552339127976335
$ws.Range("J5").Value = 0.1552339127976336
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 212.29926542574
$ws.Range("R5").Value = 1910.69338883166
$ws.Range("S5").Value = 0.07278644642688128
$ws.Range("T5").Value = 0.0727864464268813

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.428335
$ws.Range("H6").Value = 88.285005
$ws.Range("I6").Value = 0.1552339127976335
$ws.Range("J6").Value = 0.1552339127976336
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 209.24578138614
$ws.Range("R6").Value = 1883.21203247526
$ws.Range("S6").Value = 0.07173956455464314
$ws.Range("T6").Value = 0.07173956455464316

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.428335
$ws.Range("H7").Value = 88.285005
$ws.Range("I7").Value = 0.1552339127976335
$ws.Range("J7").Value = 0.1552339127976336
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 31.23218403160499
$ws.Range("R7").Value = 281.089656284445
$ws.Range("S7").Value = 0.01070790181610912
$ws.Range("T7").Value = 0.01070790181610913

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 145.9323983333333
$ws.Range("H8").Value = 437.797195
$ws.Range("I8").Value = 0.7697906523500631
$ws.Range("J8").Value = 0.7697906523500631
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 1052.772471428749
$ws.Range("R8").Value = 9474.952242858739
$ws.Range("S8").Value = 0.3609412728662858
$ws.Range("T8").Value = 0.3609412728662858

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 145.9323983333333
$ws.Range("H9").Value = 437.797195
$ws.Range("I9").Value = 0.7697906523500631
$ws.Range("J9").Value = 0.7697906523500631
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 1037.630525777682
$ws.Range("R9").Value = 9338.674731999139
$ws.Range("S9").Value = 0.3557498822426775
$ws.Range("T9").Value = 0.3557498822426775

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 145.9323983333333
$ws.Range("H10").Value = 437.797195
$ws.Range("I10").Value = 0.7697906523500631
$ws.Range("J10").Value = 0.7697906523500631
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 154.8775192657061
$ws.Range("R10").Value = 1393.897673391355
$ws.Range("S10").Value = 0.05309949724109979
$ws.Range("T10").Value = 0.05309949724109979
